$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B1").Value = "Children Names"
$ws.Range("C1").Value = "Compleat Time(s)"
$ws.Range("D1").Value = "Punctuation"
$ws.Range("G1").Value = "Words"
$ws.Range("H1").Value = "Fails"
$ws.Range("J1").Value = "Names"
$ws.Range("K1").Value = "Words"
$ws.Range("L1").Value = "Fails"

# Children data table (B:D)
$ws.Range("B2").Value = "Paco"
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 67

$ws.Range("B3").Value = "Pepe"
$ws.Range("C3").Value = 23
$ws.Range("D3").Value = 35

$ws.Range("B4").Value = "Bea"
$ws.Range("C4").Value = 54
$ws.Range("D4").Value = 24

$ws.Range("B5").Value = "Rubén"
$ws.Range("C5").Value = 32
$ws.Range("D5").Value = 47

$ws.Range("B6").Value = "María"
$ws.Range("C6").Value = 122
$ws.Range("D6").Value = 33

# Words / Fails summary table (G:H)
$ws.Range("G2").Value = "galdiolo"
$ws.Range("H2").Value = 12

$ws.Range("G3").Value = "flor"
$ws.Range("H3").Value = 4

$ws.Range("G4").Value = "palmera"
$ws.Range("H4").Value = 16

$ws.Range("G5").Value = "bloso"
$ws.Range("H5").Value = 3

# Per-child word fail tables (J:L)
# Paco
$ws.Range("J2").Value = "Paco:"
$ws.Range("K2").Value = "galdiolo"
$ws.Range("L2").Value = 1
$ws.Range("K3").Value = "flor"
$ws.Range("L3").Value = 0
$ws.Range("K4").Value = "palmera"
$ws.Range("L4").Value = 0
$ws.Range("K5").Value = "bloso"
$ws.Range("L5").Value = 1

# Pepe
$ws.Range("J7").Value = "Pepe:"
$ws.Range("K7").Value = "galdiolo"
$ws.Range("L7").Value = 1
$ws.Range("K8").Value = "flor"
$ws.Range("L8").Value = 0
$ws.Range("K9").Value = "palmera"
$ws.Range("L9").Value = 0
$ws.Range("K10").Value = "bloso"
$ws.Range("L10").Value = 1

# Bea
$ws.Range("J12").Value = "Bea:"
$ws.Range("K12").Value = "galdiolo"
$ws.Range("L12").Value = 1
$ws.Range("K13").Value = "flor"
$ws.Range("L13").Value = 1
$ws.Range("K14").Value = "palmera"
$ws.Range("L14").Value = 1
$ws.Range("K15").Value = "bloso"
$ws.Range("L15").Value = 1

# Rubén
$ws.Range("J17").Value = "Rubén:"
$ws.Range("K17").Value = "galdiolo"
$ws.Range("L17").Value = 0
$ws.Range("K18").Value = "flor"
$ws.Range("L18").Value = 0
$ws.Range("K19").Value = "palmera"
$ws.Range("L19").Value = 0
$ws.Range("K20").Value = "bloso"
$ws.Range("L20").Value = 1

# María
$ws.Range("J22").Value = "María:"
$ws.Range("K22").Value = "galdiolo"
$ws.Range("L22").Value = 0
$ws.Range("K23").Value = "flor"
$ws.Range("L23").Value = 1
$ws.Range("K24").Value = "palmera"
$ws.Range("L24").Value = 0
$ws.Range("K25").Value = "bloso"
$ws.Range("L25").Value = 1
